# "Allow Bonding with client" - adds Security Manager bonding info to the
# Server Command Table (A8), plus a note about the PB0 passkey-confirm flow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 (under the "Events only for Slaves/Server" / security-manager
# related rows): commands to call, and a short note label.
$ws.Range("E11").Value = "sl_bt_sm_configure()" + [char]10 + "sl_bt_sm_set_bondable_mode()"
$ws.Range("F11").Value = "Security  manager, "

# Row 11 needs to grow to fit the two-line command text.
$ws.Rows.Item(11).RowHeight = 37.5

# Row 28: note on how the passkey confirmation is actually triggered
# (PB0 press). Pick up the same (un-tinted) font styling already used by
# the similar note in E24, rather than the teal-ish default for this column.
$ws.Range("E24").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = "waitForPB0Press()" + [char]10 + "sl_bt_sm_passkey_confirm()"

# Restore the selection to where the user left off editing.
[void]$ws.Range("E29").Select()
